$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.982.92"
$ws.Range("E2").Value = "'  -1.10%  "
$ws.Range("D3").Value = "'2.043.83"
$ws.Range("E3").Value = "'  -1.70%  "
$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = "'  +0.39%  "
$ws.Range("D5").Value = "'250.57"
$ws.Range("E5").Value = "'  -0.25%  "
$ws.Range("E6").Value = "'  +0.08%  "
$ws.Range("E7").Value = "'  -0.01%  "
$ws.Range("D8").Value = "'56.28"
$ws.Range("E8").Value = "'  -1.00%  "
$ws.Range("D9").Value = "'61.26"
$ws.Range("E9").Value = "'  -1.53%  "
$ws.Range("E10").Value = "'  -1.01%  "
$ws.Range("E11").Value = "'  +3.07%  "
$ws.Range("D13").Value = "'16.42"
$ws.Range("E13").Value = "'  +5.09%  "
$ws.Range("D14").Value = "'2.342.10"
$ws.Range("E14").Value = "'  -1.67%  "
$ws.Range("D15").Value = "'0.798"
$ws.Range("E15").Value = "'  -5.57%  "
$ws.Range("D16").Value = "'5.53"
$ws.Range("E16").Value = "'  +4.18%  "
$ws.Range("D17").Value = "'2.038.10"
$ws.Range("E17").Value = "'  -1.99%  "
$ws.Range("D18").Value = "'36.915.32"
$ws.Range("E18").Value = "'  -1.02%  "
$ws.Range("E19").Value = "'  +12.65%  "
$ws.Range("D20").Value = "'74.86"
$ws.Range("E20").Value = "'  +2.19%  "
$ws.Range("D21").Value = "'0.0₃0902"
$ws.Range("E21").Value = "'  +5.92%  "
$ws.Range("D22").Value = "'5.36"
$ws.Range("E22").Value = "'  +1.32%  "
$ws.Range("D23").Value = "'236.23"
$ws.Range("E23").Value = "'  -2.05%  "
$ws.Range("E24").Value = "'  +0.18%  "
$ws.Range("E25").Value = "'  -3.97%  "
$ws.Range("D26").Value = "'2.37"
$ws.Range("E26").Value = "'  +17.20%  "
$ws.Range("D27").Value = "'168.88"
$ws.Range("E27").Value = "'  -1.68%  "
$ws.Range("D28").Value = "'9.22"
$ws.Range("E28").Value = "'  -0.16%  "
$ws.Range("D29").Value = "'20.07"
$ws.Range("E29").Value = "'  -4.58%  "
$ws.Range("E30").Value = "'  +0.42%  "
$ws.Range("D31").Value = "'1.17"
$ws.Range("E31").Value = "'  +5.55%  "
$ws.Range("E32").Value = "'  +2.08%  "
$ws.Range("D33").Value = "'0.0617"
$ws.Range("E33").Value = "'  -2.16%  "
$ws.Range("E34").Value = "'  +1.47%  "
$ws.Range("D35").Value = "'0.0887"
$ws.Range("E35").Value = "'  -2.39%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "'  +0.10%  "
$ws.Range("E37").Value = "'  -3.62%  "
$ws.Range("E38").Value = "'  -4.75%  "
$ws.Range("E39").Value = "'  +8.97%  "
$ws.Range("E40").Value = "'  -1.02%  "
$ws.Range("D41").Value = "'17.57"
$ws.Range("E41").Value = "'  -1.15%  "
$ws.Range("D42").Value = "'0.0222"
$ws.Range("E42").Value = "'  -2.93%  "
$ws.Range("E43").Value = "'  -3.58%  "
$ws.Range("D44").Value = "'96.38"
$ws.Range("E44").Value = "'  -3.78%  "
$ws.Range("E45").Value = "'  +1.29%  "
$ws.Range("D46").Value = "'4.71"
$ws.Range("E46").Value = "'  +16.51%  "
$ws.Range("E47").Value = "'  +2.83%  "
$ws.Range("D48").Value = "'1.279.65"
$ws.Range("E48").Value = "'  -3.65%  "
$ws.Range("E49").Value = "'  -1.41%  "
$ws.Range("E50").Value = "'  -4.56%  "
$ws.Range("D51").Value = "'2.225.97"
$ws.Range("E51").Value = "'  -1.50%  "
